$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add three new product rows (98, 99, 100) mirroring existing layout:
# A=urun_adi, B=fiyat, C=kategori, D=gorsel, E=aciklama, F=stok

$ws.Range("A98").Value = "ERKEK KAPİTONE CEKET"
$ws.Range("B98").Value = "400 TL"
$ws.Range("C98").Value = "Ceket"
$ws.Range("D98").Value = "kapitone3.jpg"
$ws.Range("E98").Value = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F98").Value = "Var"

$ws.Range("A99").Value = "ERKEK KAPİTONE CEKET"
$ws.Range("B99").Value = "400 TL"
$ws.Range("C99").Value = "Ceket"
$ws.Range("D99").Value = "kapitone2.jpg"
$ws.Range("E99").Value = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F99").Value = "Var"

$ws.Range("A100").Value = "ERKEK KAPİTONE CEKET"
$ws.Range("B100").Value = "400 TL"
$ws.Range("C100").Value = "Ceket"
$ws.Range("D100").Value = "kapitone1.jpg"
$ws.Range("E100").Value = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F100").Value = "Var"

$ws.Range("E99:E100").Select()

$wb.Save()
